# Update column F (dSF) values on Sheet1 per repulled data.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 3
    3  = -1
    4  = -1
    5  = 4
    6  = -2
    7  = 1
    8  = 2
    9  = 2
    10 = -1
    11 = -2
    13 = 5
    14 = -1
    15 = -8
    16 = 1
    18 = 5
    19 = 11
    20 = 6
    22 = 5
    23 = 1
    24 = -2
    25 = -1
    26 = 4
    28 = 2
    29 = -5
    30 = 7
    31 = 9
    32 = 9
    33 = -2
    34 = -2
    36 = 1
}

foreach ($row in $values.Keys) {
    $ws.Cells.Item($row, 6).Value = $values[$row]
}
